$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns after CS1: "Gender Identity" (CT1) and
# "Sexual Orientation" (CU1), matching the style of the preceding header
# cell (CS1) so they blend in with the rest of the header row.
$ws.Range("CS1").Copy() | Out-Null
$ws.Range("CT1:CU1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("CT1").Value = "Gender Identity"
$ws.Range("CU1").Value = "Sexual Orientation"

# Move the active selection, matching the author's final cursor position.
$ws.Range("CN9").Select() | Out-Null
